# Review protocol update: close out findings rows 3-9 with a "Closed on"
# date, and resolve the "Smoothness of resized images is not good" finding
# (row 9) with a status + remark describing the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Findings")

# Widen the "Closed on" column slightly now that it holds dates (target
# 17.140625 chars; the host's pixel-snapped ColumnWidth setter lands on the
# nearest representable width, ~17.17, from this input).
$ws.Columns.Item(6).ColumnWidth = 16.3

# Rows 3-9: stamp column F ("Closed on") with 2019-09-08 (serial 43716),
# matching the same date number format already used by column A.
$dateFormat = "[`$-1010409]d\ mmmm\ yyyy;@"
foreach ($r in 3..9) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.NumberFormat = $dateFormat
    $cell.Value = 43716
}

# Row 9: the INTER_AREA fix closes this finding — flip status to Closed
# and add the remark describing the solution.
$ws.Cells.Item(9, 5).Value = "Closed"
$ws.Cells.Item(9, 7).Value = "Solution: change Sampling/Interpolation method to INTER_AREA"

# Extend the "date between 2019-08-18 and 2020-03-23" validation rule to the
# newly-filled F3:F9 cells (same rule already applied to A2:A100).
$ws.Range("F3:F9").Validation.Add(4, 1, 1, "43695", "43830")

# Move the active selection to F3 to reflect where the edit was made.
$ws.Range("F3").Select()
